# "finalizes the one-file RPSConsole app with classes"
#
# On slide 4 ("RPS OOP Planning") there is a planning table. Two cells in
# the 2nd column (rows "Gender" and "Wins") were empty placeholders and get
# filled in with the player markers "P1" / "p2", written in red, to note
# which class/column holds player-1 vs player-2 state.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$tbl = $s.Shapes.Item(2).Table

# Row 5 = "Gender" row, column 2 -> "P1" (red)
$cell1 = $tbl.Cell(5, 2).Shape.TextFrame.TextRange
$cell1.Text = "P1"
$cell1.Font.Color.RGB = 255   # RGB(255,0,0) -> red

# Row 6 = "Wins" row, column 2 -> "p2" (red)
$cell2 = $tbl.Cell(6, 2).Shape.TextFrame.TextRange
$cell2.Text = "p2"
$cell2.Font.Color.RGB = 255   # RGB(255,0,0) -> red
